$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 stays the same: "Create and Delete CitizenShip From Excel", "PASSED", "chrome"

# Rows 2-6: "States testing with JDBC", "FAILED", "chrome"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = "States testing with JDBC"
    $ws.Cells.Item($r, 2).Value = "FAILED"
    $ws.Cells.Item($r, 3).Value = "chrome"
}

# Row 7: "States testing with JDBC", "PASSED", "chrome"
$ws.Cells.Item(7, 1).Value = "States testing with JDBC"
$ws.Cells.Item(7, 2).Value = "PASSED"
$ws.Cells.Item(7, 3).Value = "chrome"
